$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column B's new cell to stay a literal text date string instead of
# being auto-coerced into an Excel date serial number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("A3").Value = 'Portfolio_2'
$ws.Range("B3").Value = '2024-09-15'
$ws.Range("C3").Value = '[{''id'': 32981, ''symbol'': ''CATE'', ''name'': ''Cate'', ''price'': 0.04622151630505624, ''change_24h'': 15.65927489, ''change_7d'': 525.26505339, ''market_cap'': 46221516.305056244}, {''id'': 32902, ''symbol'': ''BRUH'', ''name'': ''Bruh'', ''price'': 0.00031770425887260134, ''change_24h'': -11.81050116, ''change_7d'': -54.21449484, ''market_cap'': 317704.2588726013}, {''id'': 32862, ''symbol'': ''MUNCAT'', ''name'': ''MUNCAT'', ''price'': 0.020332950047081428, ''change_24h'': 0.04479229, ''change_7d'': 24.47954696, ''market_cap'': 20332950.04708143}, {''id'': 32844, ''symbol'': ''WIWI'', ''name'': ''Wiggly Willy'', ''price'': 1.2794438148911169e-06, ''change_24h'': 0.36519617, ''change_7d'': -25.90473384, ''market_cap'': 832663.1003376707}, {''id'': 32698, ''symbol'': ''DOGS'', ''name'': ''DOGS'', ''price'': 0.001070001016141869, ''change_24h'': -2.1870609, ''change_7d'': 5.78861839, ''market_cap'': 552923025.0913107}, {''id'': 32821, ''symbol'': ''BULLS'', ''name'': ''BULLS'', ''price'': 0.0003010673566050869, ''change_24h'': -0.31866993, ''change_7d'': -72.53269195, ''market_cap'': 301067.3566050869}, {''id'': 32797, ''symbol'': ''IVfun'', ''name'': ''Invest Zone'', ''price'': 0.018351533177714856, ''change_24h'': -5.90489481, ''change_7d'': -33.83488965, ''market_cap'': 17254458.447773766}, {''id'': 32790, ''symbol'': ''EBULL'', ''name'': ''ETHEREUM IS GOOD'', ''price'': 0.0002474804052799423, ''change_24h'': 0.9636044, ''change_7d'': -45.5373862, ''market_cap'': 2474804.052799423}, {''id'': 32732, ''symbol'': ''SUNPEPE'', ''name'': ''sunpepe'', ''price'': 0.000897168648582681, ''change_24h'': 1.58911547, ''change_7d'': -23.06688199, ''market_cap'': 897168.648582681}, {''id'': 32727, ''symbol'': ''FOFAR'', ''name'': ''Fofar (fofar.meme)'', ''price'': 0.005994983027720214, ''change_24h'': 0.5843079, ''change_7d'': 2.52146706, ''market_cap'': 5695188.11063377}, {''id'': 32724, ''symbol'': ''CAT'', ''name'': "Simon''s Cat", ''price'': 3.500854798810511e-05, ''change_24h'': 5.00455684, ''change_7d'': 21.68004737, ''market_cap'': 236306135.10035145}, {''id'': 32717, ''symbol'': ''SUNDOG'', ''name'': ''SUNDOG'', ''price'': 0.34380402988132813, ''change_24h'': 0.42014554, ''change_7d'': 41.40471499, ''market_cap'': 342943334.3703298}, {''id'': 32638, ''symbol'': ''CATDOG'', ''name'': ''CatDog'', ''price'': 7.402990858692653e-05, ''change_24h'': -10.65482375, ''change_7d'': 3.12557669, ''market_cap'': 7402990.858692653}, {''id'': 32625, ''symbol'': ''FUKU'', ''name'': ''Fuku-Kun'', ''price'': 1.5877577345751996e-06, ''change_24h'': -9.703381, ''change_7d'': -47.90988549, ''market_cap'': 510971.83034441207}, {''id'': 32618, ''symbol'': ''WDOG'', ''name'': ''Wrapped Dog'', ''price'': 0.013227097510162796, ''change_24h'': -10.1036857, ''change_7d'': -24.8197901, ''market_cap'': 13227076.862663582}, {''id'': 32598, ''symbol'': ''GODCAT'', ''name'': ''Godcat Exploding Kittens'', ''price'': 0.0001428806379538942, ''change_24h'': -19.76930028, ''change_7d'': 76.51318518, ''market_cap'': 113865.6434821598}, {''id'': 32590, ''symbol'': ''DOGY'', ''name'': ''Dogy'', ''price'': 1.1893418755613293e-07, ''change_24h'': -0.5283509, ''change_7d'': -29.43293235, ''market_cap'': 103457.79178549594}, {''id'': 32524, ''symbol'': ''GINNAN'', ''name'': ''Ginnan The Cat'', ''price'': 2.9331244838977883e-06, ''change_24h'': -10.2501768, ''change_7d'': -15.88496351, ''market_cap'': 20238554.134991195}, {''id'': 32521, ''symbol'': ''NEIRO'', ''name'': ''First Neiro On Ethereum'', ''price'': 4.208864094761428e-05, ''change_24h'': 0.02630109, ''change_7d'': -25.71518583, ''market_cap'': 17411948.878136143}, {''id'': 32495, ''symbol'': ''$YAWN'', ''name'': ''YAWN'', ''price'': 3.691291491467053e-05, ''change_24h'': -1.13594427, ''change_7d'': -8.88183076, ''market_cap'': 236816.17226191767}, {''id'': 32464, ''symbol'': ''NEIRO'', ''name'': ''Neiro (neiro.lol)'', ''price'': 0.003357283483090154, ''change_24h'': -9.83454865, ''change_7d'': -24.45073389, ''market_cap'': 3357283.379014366}, {''id'': 32462, ''symbol'': ''YOUNES'', ''name'': ''YOUNES'', ''price'': 0.0019201691760439004, ''change_24h'': 0.04371067, ''change_7d'': 4.47551325, ''market_cap'': 1920165.7216595528}, {''id'': 32461, ''symbol'': ''NEIRO'', ''name'': ''Neiro Ethereum'', ''price'': 0.13394809708361538, ''change_24h'': -19.45958405, ''change_7d'': -8.80093972, ''market_cap'': 133948097.08361538}, {''id'': 32448, ''symbol'': ''BRAINLET'', ''name'': ''Brainlet'', ''price'': 0.015264458273555239, ''change_24h'': -11.10895338, ''change_7d'': -22.56562209, ''market_cap'': 14545586.869234081}, {''id'': 32415, ''symbol'': ''SKBDI'', ''name'': ''Skibidi Toilet'', ''price'': 0.05019684172868577, ''change_24h'': -3.35381365, ''change_7d'': -18.52161743, ''market_cap'': 3515720.1332713356}, {''id'': 32350, ''symbol'': ''SOY'', ''name'': ''Soyjak'', ''price'': 0.0006979898173231251, ''change_24h'': -12.64424242, ''change_7d'': -39.59443819, ''market_cap'': 646728.0414813921}, {''id'': 32310, ''symbol'': ''SAITAMA'', ''name'': ''SAITAMA INU'', ''price'': 0.00033670238715070605, ''change_24h'': -4.92073736, ''change_7d'': 51.23960976, ''market_cap'': 224403.18071201874}, {''id'': 32297, ''symbol'': ''BILL'', ''name'': ''BILL THE BEAR'', ''price'': 2.105708232664283e-06, ''change_24h'': -10.0729372, ''change_7d'': 2.45341264, ''market_cap'': 1871734.1356051096}, {''id'': 32288, ''symbol'': ''torsy'', ''name'': ''TORSY'', ''price'': 0.003067949816923183, ''change_24h'': -1.07737468, ''change_7d'': 9.19390548, ''market_cap'': 2147563.868626638}, {''id'': 32259, ''symbol'': ''MOBY'', ''name'': ''Moby'', ''price'': 0.00014163620888763475, ''change_24h'': -6.02447285, ''change_7d'': 20.19311015, ''market_cap'': 1416362.0888763475}, {''id'': 32128, ''symbol'': ''MSI'', ''name'': ''monkey shit inu'', ''price'': 1.6455146724593568e-06, ''change_24h'': 17.82578479, ''change_7d'': 30.23824405, ''market_cap'': 4844839.220083151}, {''id'': 32125, ''symbol'': ''COIN'', ''name'': ''8-Bit Coin'', ''price'': 0.0006392985397571367, ''change_24h'': -3.72990149, ''change_7d'': -16.293883, ''market_cap'': 19178956.192714103}, {''id'': 32074, ''symbol'': ''FEG'', ''name'': ''FEED EVERY GORILLA'', ''price'': 0.00032864666641601437, ''change_24h'': -6.86729081, ''change_7d'': -14.19314604, ''market_cap'': 26469704.97511844}, {''id'': 32019, ''symbol'': ''CRASH'', ''name'': ''Crash'', ''price'': 0.0035908137716047415, ''change_24h'': 4.74084861, ''change_7d'': 6.93844363, ''market_cap'': 3512668.1015975485}, {''id'': 31923, ''symbol'': ''HAWKTUAH'', ''name'': ''Hawk Tuah'', ''price'': 0.0008897625533380383, ''change_24h'': -6.65733615, ''change_7d'': 10.34980199, ''market_cap'': 764056.5777083959}, {''id'': 31921, ''symbol'': ''RETARDIO'', ''name'': ''RETARDIO'', ''price'': 0.040150477657181494, ''change_24h'': 2.87572677, ''change_7d'': 7.38531031, ''market_cap'': 39173667.32613651}, {''id'': 31914, ''symbol'': ''BILLY'', ''name'': ''Billy'', ''price'': 0.02036048834136512, ''change_24h'': -12.99161381, ''change_7d'': -12.71696274, ''market_cap'': 19060219.85126136}, {''id'': 31908, ''symbol'': ''WAT'', ''name'': ''Wat'', ''price'': 7.900638295374549e-06, ''change_24h'': -3.89456196, ''change_7d'': 34.25522638, ''market_cap'': 3305511.7276330935}, {''id'': 31847, ''symbol'': ''WOLF'', ''name'': ''LandWolf (SOL)'', ''price'': 0.002000247992751642, ''change_24h'': -5.10306212, ''change_7d'': 27.89328079, ''market_cap'': 20002166.270074576}, {''id'': 31830, ''symbol'': ''DADDY'', ''name'': ''Daddy Tate'', ''price'': 0.07541035601469624, ''change_24h'': 2.3400481, ''change_7d'': 0.79026391, ''market_cap'': 45222710.38962473}, {''id'': 31798, ''symbol'': ''JENNER'', ''name'': ''Caitlyn Jenner (ETH)'', ''price'': 0.0002215635197330867, ''change_24h'': 3.22030906, ''change_7d'': 3.81997045, ''market_cap'': 213714.15834373693}, {''id'': 31770, ''symbol'': ''HONK'', ''name'': ''Pepoclown'', ''price'': 2.474355501280837e-09, ''change_24h'': -9.45214316, ''change_7d'': 85.55821414, ''market_cap'': 992845.525569602}, {''id'': 31678, ''symbol'': ''MAGAA'', ''name'': ''MAGA AGAIN'', ''price'': 0.007149135149749142, ''change_24h'': -7.00163058, ''change_7d'': -6.05357947, ''market_cap'': 7149135.149749142}, {''id'': 31632, ''symbol'': ''PEIPEI'', ''name'': ''PeiPei (ETH)'', ''price'': 1.1149213635392922e-07, ''change_24h'': -1.19243483, ''change_7d'': -12.35441257, ''market_cap'': 45840699.13043585}, {''id'': 31569, ''symbol'': ''TROG'', ''name'': ''Trog'', ''price'': 2.249779085048139e-06, ''change_24h'': 1.43307773, ''change_7d'': 48.51745419, ''market_cap'': 946459.5632889016}, {''id'': 31561, ''symbol'': ''SOLCAT'', ''name'': ''CatSolHat'', ''price'': 0.012559164084184048, ''change_24h'': -0.17625545, ''change_7d'': 6.068356, ''market_cap'': 830901.2845156254}, {''id'': 31510, ''symbol'': ''MOTHER'', ''name'': ''Mother Iggy'', ''price'': 0.05950302362871237, ''change_24h'': 9.13057155, ''change_7d'': 66.53466622, ''market_cap'': 58678499.412200876}, {''id'': 31496, ''symbol'': ''PAPU'', ''name'': ''Papu Token'', ''price'': 8.127591752575577e-09, ''change_24h'': 3.32237823, ''change_7d'': 84.29785943, ''market_cap'': 257687.8581491088}, {''id'': 31317, ''symbol'': ''MEOW'', ''name'': ''MeowCat'', ''price'': 0.008027174921922751, ''change_24h'': -1.93791851, ''change_7d'': -4.6879953, ''market_cap'': 802717.4921922751}, {''id'': 31305, ''symbol'': ''MAGA'', ''name'': ''MAGA'', ''price'': 6.45821541829669e-05, ''change_24h'': 2.61448697, ''change_7d'': -33.55051374, ''market_cap'': 25203760.182413302}, {''id'': 31284, ''symbol'': ''HAMMY'', ''name'': ''SAD HAMSTER'', ''price'': 0.012036529348068124, ''change_24h'': -7.93810938, ''change_7d'': 16.15741239, ''market_cap'': 11337222.038607052}, {''id'': 31267, ''symbol'': ''UTYAB'', ''name'': ''UTYABSWAP'', ''price'': 0.0013080123797746539, ''change_24h'': 12.87982157, ''change_7d'': 90.98739274, ''market_cap'': 1308012.3797746538}, {''id'': 31259, ''symbol'': ''FOXSY'', ''name'': ''Foxsy AI'', ''price'': 0.09002360660849407, ''change_24h'': -1.0132708, ''change_7d'': 15.70748931, ''market_cap'': 18296466.67516939}, {''id'': 31163, ''symbol'': ''SLOTH'', ''name'': ''Slothana'', ''price'': 0.006046333945347184, ''change_24h'': -2.7642754, ''change_7d'': 2.40574951, ''market_cap'': 10249250.423811788}, {''id'': 31152, ''symbol'': ''KENDU'', ''name'': ''Kendu Inu'', ''price'': 5.4935096559914444e-05, ''change_24h'': -0.97786411, ''change_7d'': 20.00854328, ''market_cap'': 52101991.84975898}, {''id'': 31121, ''symbol'': ''LABZ'', ''name'': ''Insane Labz'', ''price'': 0.01918521121245546, ''change_24h'': -3.11082122, ''change_7d'': 3.79800696, ''market_cap'': 949806.9599914897}, {''id'': 31106, ''symbol'': ''DOKY'', ''name'': ''Donkey King'', ''price'': 2.2075139455454384e-05, ''change_24h'': -5.01876988, ''change_7d'': -28.3249804, ''market_cap'': 191492.16518065598}, {''id'': 31061, ''symbol'': ''DUREV'', ''name'': ''Povel Durev'', ''price'': 0.035918555440704596, ''change_24h'': -0.4457259, ''change_7d'': 2.74063729, ''market_cap'': 3412262.766866937}, {''id'': 31051, ''symbol'': ''LOBO'', ''name'': ''LOBO•THE•WOLF•PUP'', ''price'': 0.0005446515319582705, ''change_24h'': -4.39495922, ''change_7d'': -14.72250201, ''market_cap'': 8925630.381070133}, {''id'': 31044, ''symbol'': ''HEGE'', ''name'': ''Hege'', ''price'': 0.015461313982300686, ''change_24h'': 1.49023661, ''change_7d'': 17.33684212, ''market_cap'': 15459031.843962986}, {''id'': 31036, ''symbol'': ''CARLO'', ''name'': ''Carlo'', ''price'': 0.0009803795167569037, ''change_24h'': -4.01065896, ''change_7d'': -12.71171748, ''market_cap'': 908744.9658366792}, {''id'': 30979, ''symbol'': ''ELON'', ''name'': ''Elon MemeLord'', ''price'': 4.842346651796097e-06, ''change_24h'': 13.13773278, ''change_7d'': 20.27182781, ''market_cap'': 3318071.7158579836}, {''id'': 30968, ''symbol'': ''LONG'', ''name'': ''Long'', ''price'': 1.1058038310879761e-06, ''change_24h'': -0.43797833, ''change_7d'': 32.47623643, ''market_cap'': 866633.7614739848}, {''id'': 30953, ''symbol'': ''SHIBTC'', ''name'': ''ShibaBitcoin'', ''price'': 0.6975613166927364, ''change_24h'': -0.09149543, ''change_7d'': -0.10271175, ''market_cap'': 60311902.01723075}, {''id'': 30943, ''symbol'': ''$MICHI'', ''name'': ''michi (SOL)'', ''price'': 0.117911869351425, ''change_24h'': -2.72345559, ''change_7d'': 14.56403438, ''market_cap'': 65531812.063366674}, {''id'': 30933, ''symbol'': ''DOG'', ''name'': ''Dog (Runes)'', ''price'': 0.0021882683829385547, ''change_24h'': -7.30841919, ''change_7d'': -13.16273772, ''market_cap'': 218826838.29385546}, {''id'': 30912, ''symbol'': ''MANEKI'', ''name'': ''MANEKI'', ''price'': 0.005424160880608815, ''change_24h'': -1.55428908, ''change_7d'': 27.75077777, ''market_cap'': 45687095.31711064}, {''id'': 30867, ''symbol'': ''WHY'', ''name'': ''WHY'', ''price'': 2.678935018639163e-07, ''change_24h'': -2.09424985, ''change_7d'': -1.6444193, ''market_cap'': 112515270.78284486}, {''id'': 30859, ''symbol'': ''HOPPY'', ''name'': ''Hoppy'', ''price'': 5.7326445931293125e-05, ''change_24h'': -8.5522925, ''change_7d'': 19.82919212, ''market_cap'': 24116662.538835704}, {''id'': 30828, ''symbol'': ''CATA'', ''name'': ''CATAMOTO'', ''price'': 0.00027747981086225166, ''change_24h'': -1.00979801, ''change_7d'': 19.74747714, ''market_cap'': 5376911.293642551}, {''id'': 30647, ''symbol'': ''GOAT'', ''name'': ''Sonic The Goat'', ''price'': 3.2588968825515676e-07, ''change_24h'': 0.33699756, ''change_7d'': -5.31891005, ''market_cap'': 288771.7766670043}, {''id'': 30645, ''symbol'': ''POWSCHE'', ''name'': ''Powsche'', ''price'': 0.017237224992953708, ''change_24h'': -1.21393757, ''change_7d'': 30.59615048, ''market_cap'': 1638391.9267559026}, {''id'': 30629, ''symbol'': ''KITTENWIF'', ''name'': ''KittenWifHat'', ''price'': 0.0008244819696452965, ''change_24h'': 0.55046486, ''change_7d'': 3.14689902, ''market_cap'': 805996.554023717}, {''id'': 30601, ''symbol'': ''MOUTAI'', ''name'': ''Moutai'', ''price'': 0.006432758184363277, ''change_24h'': -6.29420006, ''change_7d'': -6.73813125, ''market_cap'': 4481814.836773734}, {''id'': 30550, ''symbol'': ''HAMI'', ''name'': ''HAMI'', ''price'': 0.001750194661509419, ''change_24h'': 2.40717766, ''change_7d'': 8.38800081, ''market_cap'': 1684101.8947695398}, {''id'': 30493, ''symbol'': ''NUB'', ''name'': ''nubcat'', ''price'': 0.008602029181940361, ''change_24h'': -5.96381558, ''change_7d'': -13.65462854, ''market_cap'': 8171929.899156726}, {''id'': 30484, ''symbol'': ''PUPS'', ''name'': ''PUPS (Ordinals)'', ''price'': 2.719228707186724, ''change_24h'': 0.37133471, ''change_7d'': -6.98688497, ''market_cap'': 21128407.054840844}, {''id'': 30409, ''symbol'': ''NSO'', ''name'': "NeverSurrenderOne''s", ''price'': 0.00019396208554128963, ''change_24h'': -2.26618058, ''change_7d'': 3.35814086, ''market_cap'': 124758.63040679523}, {''id'': 30407, ''symbol'': ''ROOST'', ''name'': ''Roost Coin'', ''price'': 0.0020775869321590043, ''change_24h'': 16.27924256, ''change_7d'': 50.55261981, ''market_cap'': 1946997.6521838105}, {''id'': 30402, ''symbol'': ''CAW'', ''name'': ''crow with knife'', ''price'': 3.766828728058379e-08, ''change_24h'': -3.21365199, ''change_7d'': 4.23995501, ''market_cap'': 28999360.29598486}, {''id'': 30361, ''symbol'': ''DOGEMOB'', ''name'': ''DOGEMOB'', ''price'': 0.00034782268006846326, ''change_24h'': -4.12544203, ''change_7d'': -1.35392858, ''market_cap'': 173720.06955988053}, {''id'': 30309, ''symbol'': ''SC'', ''name'': ''Shark Cat'', ''price'': 0.012052976711005805, ''change_24h'': -6.4461373, ''change_7d'': -11.71048065, ''market_cap'': 11931187.636836004}, {''id'': 30303, ''symbol'': ''GB'', ''name'': ''Gary Banking'', ''price'': 9.672085807112003e-05, ''change_24h'': -3.01620481, ''change_7d'': -37.39824308, ''market_cap'': 1015569.0097467602}, {''id'': 30285, ''symbol'': ''MUMU'', ''name'': ''Mumu the Bull'', ''price'': 3.651115704661971e-05, ''change_24h'': -5.92629023, ''change_7d'': 2.78724668, ''market_cap'': 83362435.03807513}, {''id'': 30270, ''symbol'': ''FOMO'', ''name'': ''FOMO BULL CLUB'', ''price'': 9.464007213768345e-06, ''change_24h'': -1.05740597, ''change_7d'': -39.27013238, ''market_cap'': 316699.0380232343}, {''id'': 30218, ''symbol'': ''MABA'', ''name'': ''Make America Based Again'', ''price'': 0.00024460005932453263, ''change_24h'': 0.27982023, ''change_7d'': 14.86484953, ''market_cap'': 225620.62689572052}, {''id'': 30193, ''symbol'': ''BENJI'', ''name'': ''Basenji'', ''price'': 0.015489738445553111, ''change_24h'': -1.60147188, ''change_7d'': -13.33552972, ''market_cap'': 14282233.87136402}, {''id'': 30133, ''symbol'': ''COINYE'', ''name'': ''Coinye West'', ''price'': 0.0008938532605026727, ''change_24h'': 0.58480217, ''change_7d'': 26.76437029, ''market_cap'': 716195.2107956706}, {''id'': 30126, ''symbol'': ''MEW'', ''name'': ''cat in a dogs world'', ''price'': 0.0042752317408956565, ''change_24h'': -3.83502265, ''change_7d'': 3.98341454, ''market_cap'': 380020599.1869248}, {''id'': 30119, ''symbol'': ''WIF'', ''name'': ''DOGWIFHOOD'', ''price'': 0.00042830709499371516, ''change_24h'': -4.85299443, ''change_7d'': 9.96572314, ''market_cap'': 427844.59742824937}, {''id'': 30096, ''symbol'': ''DEGEN'', ''name'': ''Degen'', ''price'': 0.003245256240496699, ''change_24h'': -0.5219058, ''change_7d'': 8.02685185, ''market_cap'': 46016464.202377155}, {''id'': 30071, ''symbol'': ''LADYF'', ''name'': ''Milady Wif Hat'', ''price'': 1.3713102368460841e-05, ''change_24h'': -1.23880289, ''change_7d'': 3.96667031, ''market_cap'': 3428251.539550532}, {''id'': 30063, ''symbol'': ''GIGA'', ''name'': ''Gigachad'', ''price'': 0.016779323169116272, ''change_24h'': -7.15490186, ''change_7d'': 12.10071658, ''market_cap'': 156088175.32098106}, {''id'': 30011, ''symbol'': ''SKID'', ''name'': ''Success Kid'', ''price'': 0.026644795348589937, ''change_24h'': 1.68115728, ''change_7d'': 4.83207344, ''market_cap'': 2311241.436194542}, {''id'': 30008, ''symbol'': ''APU'', ''name'': ''Apu Apustaja'', ''price'': 0.00018752558223819106, ''change_24h'': -14.60615612, ''change_7d'': 3.59312395, ''market_cap'': 59698114.894198336}, {''id'': 30007, ''symbol'': ''HUND'', ''name'': ''HUND'', ''price'': 0.006624933004839104, ''change_24h'': -4.23012646, ''change_7d'': 18.80458653, ''market_cap'': 2529579.0666696564}, {''id'': 29999, ''symbol'': ''CHKN'', ''name'': ''Chickencoin'', ''price'': 5.1618896538034534e-08, ''change_24h'': 0.21755199, ''change_7d'': 2.94907457, ''market_cap'': 3315946.2947068005}, {''id'': 29932, ''symbol'': ''BSHIB'', ''name'': ''Based Shiba Inu'', ''price'': 2.0759965209992277e-05, ''change_24h'': 4.54541767, ''change_7d'': 9.29667141, ''market_cap'': 187698.96483125357}, {''id'': 29920, ''symbol'': ''SLERF'', ''name'': ''SLERF'', ''price'': 0.13567370591902544, ''change_24h'': -2.6148509, ''change_7d'': 3.48513591, ''market_cap'': 67836547.6936744}, {''id'': 29879, ''symbol'': ''ANDY'', ''name'': ''ANDY (ETH)'', ''price'': 7.275078363694599e-05, ''change_24h'': -7.99221739, ''change_7d'': 33.40693501, ''market_cap'': 72750783.636946}]'
$ws.Range("D3").Value = 59981.8
$ws.Range("E3").Value = 135.2
$ws.Range("F3").Value = '[]'

# Drop the transient "Text" number-format style so B3 ends up styled the
# same as its neighbours (no stray style index left behind).
$ws.Range("B3").Style = "Normal"
